# daily auto push: 2025-10-10 13:36 UTC
# Appends the new daily data row (row 91) to the sheet, mirroring the
# existing rows for 2025/10/10 ("金" / Friday) already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as plain text (e.g. "2025/10/10"), matching the
# other rows in the sheet. Force the cell to text format before assigning
# the value so Excel does not auto-convert the string into a date serial
# number, then restore the cell to the sheet's normal (unstyled) look so
# it matches the rest of the data rows, which carry no explicit style.
$ws.Range("A91").NumberFormat = "@"
$ws.Range("A91").Value = "2025/10/10"
$ws.Range("A91").Style = "Normal"

$ws.Range("B91").Value = "金"
$ws.Range("C91").Value = 20
$ws.Range("D91").Value = 201
